$d = $word.ActiveDocument

# Move to the very end of the document body text, then add two paragraph
# breaks (an intervening blank paragraph, matching the diff) followed by
# the new line of text - the same sequence a user would produce by
# placing the cursor at the end of the existing paragraph and pressing
# Enter twice before typing.
$sel = $word.Selection
$sel.EndKey(6) | Out-Null
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText("A change to note the git tracking")
